$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet was generically named "Sheet1"; rename it after the scraped batter.
$ws.Name = "David Miller"

# Start from a clean sheet, then rewrite the full table: a new "matchNo"
# column is inserted at the front and the single existing match row is
# joined by the rest of David Miller's scraped Rajasthan Royals innings.
$ws.Cells.Clear()

$data = @(
  @("matchNo", "teamName", "batterName", "states", "runs", "balls", "fours", "sixes", "sr", "opponentTeamName", "venue", "date", "result"),
  @("28th", "Rajasthan Royals", "David Miller", "", "7", "3", "0", "1", "233.33", "Sunrisers Hyderabad", "Delhi", "May 02", "Royals won by 55 runs"),
  @("24th", "Rajasthan Royals", "David Miller", "", "7", "4", "1", "0", "175.00", "Mumbai Indians", "Delhi", "April 29", "Mumbai won by 7 wickets (with 9 balls remaining)"),
  @("51st", "Rajasthan Royals", "David Miller", "lbw b Coulter-Nile", "15", "23", "0", "0", "65.21", "Mumbai Indians", "Sharjah", "October 05", "Mumbai won by 8 wickets (with 70 balls remaining)"),
  @("36th", "Rajasthan Royals", "David Miller", "st †Pant b Ashwin", "7", "10", "0", "0", "70.00", "Delhi Capitals", "Abu Dhabi", "September 25", "Capitals won by 33 runs"),
  @("18th", "Rajasthan Royals", "David Miller", "", "24", "23", "3", "0", "104.34", "Kolkata Knight Riders", "Wankhede", "April 24", "Royals won by 6 wickets (with 7 balls remaining)"),
  @("16th", "Rajasthan Royals", "David Miller", "lbw b Mohammed Siraj", "0", "2", "0", "0", "0.00", "Royal Challengers Bangalore", "Wankhede", "April 22", "RCB won by 10 wickets (with 21 balls remaining)"),
  @("7th", "Rajasthan Royals", "David Miller", "c Lalit Yadav b Avesh Khan", "62", "43", "7", "2", "144.18", "Delhi Capitals", "Wankhede", "April 15", "Royals won by 3 wickets (with 2 balls remaining)"),
  @("12th", "Rajasthan Royals", "David Miller", "lbw b Ali", "2", "5", "0", "0", "40.00", "Chennai Super Kings", "Wankhede", "April 19", "Super Kings won by 45 runs")
)

$quote = [string][char]39

for ($r = 0; $r -lt $data.Count; $r++) {
  $row = $data[$r]
  for ($c = 0; $c -lt $row.Count; $c++) {
    $cell = $ws.Cells.Item($r + 1, $c + 1)
    $value = $row[$c]
    if ($value -eq "") {
      # A bare "" would clear the cell to a true blank; a leading single
      # quote keeps it a (blank) text cell, matching the source data's
      # empty `states` entries for batters who were not dismissed.
      $cell.Value = $quote
    } elseif ($value -match "^[0-9]+(\.[0-9]+)?$") {
      # The scraped sheet keeps numeric-looking fields (runs, balls,
      # fours, sixes, sr) as text, so force text storage before writing
      # them or Excel would silently convert them to numbers.
      $cell.NumberFormat = "@"
      $cell.Value = $value
    } else {
      $cell.Value = $value
    }
  }
}
